$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gValues = @{
    2 = 2
    3 = 0
    4 = 2
    5 = 1
    6 = 2
    7 = 0
    8 = 0
    9 = 0
    10 = 1
    11 = 1
    12 = 2
    13 = 1
    14 = 4
    15 = 0
    16 = 1
    17 = 1
    18 = 0
    19 = 1
    20 = 2
    21 = 1
    22 = 0
    23 = 1
    24 = 2
    25 = 0
    26 = 0
    27 = 0
    29 = 0
    30 = 2
    31 = 1
    32 = 1
    33 = 2
    34 = 1
    35 = 4
    36 = 0
    37 = 3
    39 = 0
    40 = 1
    41 = 1
    42 = 0
    43 = 1
    45 = 1
    46 = 0
    47 = 0
    48 = 1
    49 = 1
    50 = 2
    51 = 3
    52 = 1
    53 = 2
    54 = 1
    55 = 2
    56 = 0
    57 = 3
    58 = 1
    59 = 1
    60 = 2
    61 = 0
    62 = 0
    63 = 0
    64 = 0
    65 = 0
    66 = 2
    67 = 0
    68 = 1
    69 = 1
    70 = 0
    71 = 2
    72 = 1
    73 = 0
    74 = 1
    75 = 0
    76 = 0
    77 = 0
    78 = 1
    79 = 1
    80 = 1
    83 = 1
    84 = 0
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}